$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 = NroPoliza, G2 = FechaSiniestro. Both cells store their values as
# text (quote-prefixed numbers/dates), so prefix the new values with an
# apostrophe to force text entry and keep the existing cell formatting
# (style / quotePrefix) untouched, matching the rest of the sheet.
$ws.Range("E2").Value = "'12112002435"
$ws.Range("G2").Value = "'23/06/2021"
